# Applies the "overview.docx" edits described by the commit diff.
$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $result = $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $result) {
        Write-Output "NOT FOUND: $old"
    }
}

# 1. Course Description paragraph
Replace-Text `
    "Every step in policymaking relies on data. This course introduces students to data management, wrangling, and visualization as well as the technical tools necessary to do such work in an open and reproducible fashion." `
    "Every step in policymaking relies on data. This course introduces students to data management, wrangling, communication, and visualization in the context of public policy, public administration, and behavioral science as well as the technical tools necessary to do such work in an open and reproducible fashion."

# 2. Expanded Description paragraph
Replace-Text `
    "Data preprocessing, wrangling, and management often consumes a large fraction of the time spent doing quantitative data analysis in public administration, public policy, and behavioral science research. Yet these topics frequently do not receive regular attention in methodological courses that focus on statistical inference. This class introduces students to the technical tools necessary to do these tasks in an open and reproducible fashion suitable for modern computational data workflows. Throughout the course of the semester, students will learn the principles and practice of conducting reproducible quantitative research, including readable programming and coding, version control, methods of documentation, data storage, workflow management, and exploratory data visualization. A variety of relevant open technical software tools will be introduced and used, including but not limited to R (and RStudio), git (and github), markdown, and a variety of helper programs to tie things together." `
    "Data preprocessing, wrangling, and management often consumes a large fraction of the time spent doing quantitative data analysis in public administration, public policy, and behavioral science research. Yet these topics frequently do not receive regular attention in methodological courses that focus on statistical inference. This class introduces students to the technical tools necessary to do these tasks in an open and reproducible fashion suitable for modern computational data workflows in the public sector. Throughout the course of the semester, students will learn the principles and practice of conducting reproducible quantitative research, including readable programming and coding, version control, methods of documentation, data storage, workflow management, and exploratory data visualization. A variety of relevant open technical software tools will be introduced and used, including but not limited to R (and RStudio), git (and github), markdown, and a variety of helper programs to tie things together. Special attention will be paid to data frequently used in public policy, public administration, and behavioral science."

# 3. Insert new "Learning Objectives" section (Heading2 + 5-item numbered list)
#    right after the Prerequisites section and before the Materials section.
$prereqIdx = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "*above all*patience.*") {
        $prereqIdx = $i
        break
    }
}
if ($prereqIdx -eq 0) {
    Write-Output "NOT FOUND: prerequisites closing paragraph"
}

$lastPrereqPara = $d.Paragraphs($prereqIdx)
$lastPrereqPara.Range.InsertParagraphAfter()

$headingPara = $d.Paragraphs($prereqIdx + 1)
$headingPara.Range.Text = "Learning Objectives"
$headingPara.Range.Style = "Heading2"

$objectives = @(
    "Demonstrate capability in open science and contemporary reproducible data analysis tools",
    "Apply appropriate principles of data and file management to data projects",
    "Evaluate the credibility and clarity of data visualizations",
    "Create effective, reproducible, and well designed data visualizations with appropriate tools",
    "Analyze large-N datasets commonly used in public policy and behavioral science"
)

$firstItemIdx = $prereqIdx + 2
$cursor = $headingPara
foreach ($objective in $objectives) {
    $cursor.Range.InsertParagraphAfter()
    $cursor = $cursor.Next()
}
for ($j = 0; $j -lt $objectives.Count; $j++) {
    $itemPara = $d.Paragraphs($firstItemIdx + $j)
    $itemPara.Range.Text = $objectives[$j]
    $itemPara.Range.Style = "Compact"
}

$lastItemIdx = $firstItemIdx + $objectives.Count - 1
$listRange = $d.Range($d.Paragraphs($firstItemIdx).Range.Start, $d.Paragraphs($lastItemIdx).Range.End)
$listRange.ListFormat.ApplyNumberDefault()

$sectionRange = $d.Range($headingPara.Range.Start, $d.Paragraphs($lastItemIdx).Range.End)
$d.Bookmarks.Add("learning-objectives", $sectionRange)

# 4. Practicum paragraph
Replace-Text `
    "The practicum is essentially a large assignment that is worth more and graded on a scale." `
    "A practicum is a large assignment that is worth more and graded on a scale."

$oldPracticum = "They are untimed, take-home, cumulative, and will be completed on your own time (and computer). Unlike the weekly assignments, you are also not allowed to work together on them. Their timing corresponds (roughly) with the midterm and final - you may wish to think of them as the " + [char]0x201c + "take-home" + [char]0x201d + " midterm and " + [char]0x201c + "take-home" + [char]0x201d + " final, if you prefer."
$newPracticum = "It is untimed, take-home, cumulative, and will be completed on your own time (and computer). Unlike the weekly assignments, you are not allowed to work together on the practicum. Essentially, think of it as take home test that complements the in-class core exam."
Replace-Text $oldPracticum $newPracticum

# 5. Core Exam paragraph
Replace-Text `
    "The core exam will be in-class. More information will be given as the exam gets closer." `
    "The core exam will be in class. More information will be given as the exam gets closer."

# 6. Final Project paragraph
Replace-Text `
    "A project utilizing data of your own choice. Graduate students will have higher expectations than undergraduate students." `
    "A project utilizing data of your own choice. Graduate students will have higher expectations than undergraduate students. More information will be given as the exam gets closer."

Write-Output "done"
